# Updates the cryptos list data (Price / Volume(1h) columns, and the
# dogwifhat/Cosmos row swap) to match the latest scrape.
#
# Setting NumberFormat to "@" (Text) before assigning .Value keeps
# numeric-looking strings (e.g. "582.48") stored as text, matching the
# source workbook's inlineStr cells instead of being coerced to
# floating-point numbers. Resetting .Style back to "Normal" afterwards
# drops the Text number-format override so the cell keeps the sheet's
# default (unstyled) appearance, same as the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.215.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.102.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.03%  "
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.89%  "
$ws.Range("E6").Style = "Normal"

# Row 7 - USDC
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E7").Style = "Normal"

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.100.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.03%  "
$ws.Range("E8").Style = "Normal"

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E9").Style = "Normal"

# Row 10 - Toncoin
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E10").Style = "Normal"

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("E11").Style = "Normal"

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.41%  "
$ws.Range("E12").Style = "Normal"

# Row 13 - ShibaInu
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("E13").Style = "Normal"

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.75%  "
$ws.Range("E14").Style = "Normal"

# Row 15 - TRON
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E15").Style = "Normal"

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.618.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.189.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("E17").Style = "Normal"

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.10%  "
$ws.Range("E18").Style = "Normal"

# Row 19 - WrappedEther
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.104.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("E19").Style = "Normal"

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +17.09%  "
$ws.Range("E20").Style = "Normal"

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("E21").Style = "Normal"

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.36%  "
$ws.Range("E22").Style = "Normal"

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("E23").Style = "Normal"

# Row 24 - Litecoin
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("E24").Style = "Normal"

# Row 25 - Fetch.AI
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.81%  "
$ws.Range("E25").Style = "Normal"

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.66%  "
$ws.Range("E26").Style = "Normal"

# Row 27 - RenderToken
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E27").Style = "Normal"

# Row 29 - NEARProtocol
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("E29").Style = "Normal"

# Row 30 - ImmutableX
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.61%  "
$ws.Range("E30").Style = "Normal"

# Row 31 - PancakeSwap
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("E31").Style = "Normal"

# Row 32 - PEPE
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("E32").Style = "Normal"

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("E33").Style = "Normal"

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.34%  "
$ws.Range("E34").Style = "Normal"

# Row 35 - FirstDigitalUSD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E35").Style = "Normal"

# Row 36 - Mantle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("E36").Style = "Normal"

# Row 37 - Filecoin
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("E37").Style = "Normal"

# Row 38 - Arweave
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.05%  "
$ws.Range("E38").Style = "Normal"

# Row 39 - Stacks
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.12%  "
$ws.Range("E39").Style = "Normal"

# Row 40 - OKB
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("E40").Style = "Normal"

# Row 41 - TheGraph
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("E41").Style = "Normal"

# Row 42 - Kaspa
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("E42").Style = "Normal"

# Row 43 - now Cosmos (was dogwifhat)
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.19%  "
$ws.Range("E43").Style = "Normal"

# Row 44 - now dogwifhat (was Cosmos)
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("E44").Style = "Normal"

# Row 45 - Bittensor
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.07%  "
$ws.Range("E45").Style = "Normal"

# Row 46 - VeChain
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("E46").Style = "Normal"

# Row 47 - Maker
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.778.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E47").Style = "Normal"

# Row 48 - Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("E48").Style = "Normal"

# Row 49 - USDe
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E49").Style = "Normal"

# Row 50 - InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.02%  "
$ws.Range("E50").Style = "Normal"

# Row 51 - ThetaToken
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.24%  "
$ws.Range("E51").Style = "Normal"
